$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 132
$ws1.Range("F3").Value = 445

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 72
$ws2.Range("F3").Value = 29

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 132
$ws4.Range("F3").Value = 72
$ws4.Range("F4").Value = 445
$ws4.Range("F8").Value = 29
